# VerveStacks_FIN / scen_tsparameters_ts_12.xlsx
# Commit: "Updated FIN model - 2025-08-28 17:12"
#
# The underlying TS-group definitions on the "ev_charging_uc" sheet
# (C13 / C14, which feed G8 / G7 via formulas) were re-shuffled, and the
# corresponding re_profiles "hydro" lookup rows (M4:N7) were re-sorted
# to match the new order.

$wb = $excel.ActiveWorkbook

# --- ev_charging_uc: update the two timeslice-group strings ----------
$wsUC = $wb.Worksheets.Item("ev_charging_uc")

# C13 feeds G8 (=C13); C14 feeds G7 (=C14) - both recalc automatically.
$wsUC.Range("C13").Value = "RaD,WaD,FaD,RaP,SaD,WaP,FaP,SaP"
$wsUC.Range("C14").Value = "RaP,FaN,SaN,WaN,FaP,SaP,RaN,WaP"

# --- re_profiles: re-order the hydro profile rows (M4:O7) -------------
$wsRE = $wb.Worksheets.Item("re_profiles")

# Row 4 <-> Row 6
$wsRE.Range("M4").Value = "R"
$wsRE.Range("N4").Value = 0.34481908618716439

$wsRE.Range("M5").Value = "S"
$wsRE.Range("N5").Value = 0.27551721102209698

$wsRE.Range("M6").Value = "F"
$wsRE.Range("N6").Value = 0.2827009419843296

$wsRE.Range("M7").Value = "W"
$wsRE.Range("N7").Value = 0.29696276080640904
